# Generate Report for Handback
#
# Row 7 in both the "zh-cn" and "de-de" sheets corresponds to the
# 41a8cd03-7af9-41c8-9ccb-4c1040ad979c file. A new handback was produced for
# it, so the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" / "Error Detail" columns (I/J/K/P) need to be
# populated (they previously held placeholder/blank values), and a
# hyperlink needs to be added on the new "Latest Target File" cell.

$wb = $excel.ActiveWorkbook

$targetFileName = "41a8cd03-7af9-41c8-9ccb-4c1040ad979c.md"
$targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d33dbdb803750b740fb47053e57a5c12430aa9e8/e2e/41a8cd03-7af9-41c8-9ccb-4c1040ad979c.md"

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("I7").Value = $targetFileName
$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $targetUrl, "", "", $targetFileName)

$wsZh.Range("J7").Value = "41a8cd03-7af9-41c8-9ccb-4c1040ad979c.59ff6f81998bce6d7e462d5f04177361aa3e2dcb.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-29 02:54:41"
$wsZh.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce8130bc931c0a39e122d5dc38dfe317219579d0/e2e/41a8cd03-7af9-41c8-9ccb-4c1040ad979c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d33dbdb803750b740fb47053e57a5c12430aa9e8/e2e/41a8cd03-7af9-41c8-9ccb-4c1040ad979c.md."

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I7").Value = $targetFileName
$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $targetUrl, "", "", $targetFileName)

$wsDe.Range("J7").Value = "41a8cd03-7af9-41c8-9ccb-4c1040ad979c.59ff6f81998bce6d7e462d5f04177361aa3e2dcb.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-29 02:54:48"
$wsDe.Range("P7").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ce8130bc931c0a39e122d5dc38dfe317219579d0/e2e/41a8cd03-7af9-41c8-9ccb-4c1040ad979c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d33dbdb803750b740fb47053e57a5c12430aa9e8/e2e/41a8cd03-7af9-41c8-9ccb-4c1040ad979c.md."
